$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.308704024249916
$ws.Range("C2").Value = 6.437078245076543
$ws.Range("D2").Value = 5.977220058499384
$ws.Range("E2").Value = 16.47857091707814
$ws.Range("G2").Value = 26.97136856850733
$ws.Range("H2").Value = 13.7802570807459
$ws.Range("K2").Value = 8.616779430938299
$ws.Range("O2").Value = 20.74820733362123

$ws.Range("B3").Value = 8.945272875311051
$ws.Range("C3").Value = 6.282184168223292
$ws.Range("D3").Value = 5.856905516016423
$ws.Range("E3").Value = 15.54580898052582
$ws.Range("G3").Value = 26.98416138509352
$ws.Range("H3").Value = 13.82916981153145
$ws.Range("K3").Value = 8.270947341183369
$ws.Range("O3").Value = 20.81612045160561

$ws.Range("B4").Value = 8.71569299596282
$ws.Range("C4").Value = 6.184375474572085
$ws.Range("D4").Value = 5.783531040022905
$ws.Range("E4").Value = 14.94833663402791
$ws.Range("G4").Value = 27.00280635262983
$ws.Range("H4").Value = 13.86180807550766
$ws.Range("K4").Value = 8.049334571137827
$ws.Range("O4").Value = 20.86322586455669

$ws.Range("B5").Value = 8.620662900412873
$ws.Range("C5").Value = 6.14387496377978
$ws.Range("D5").Value = 5.753802508808223
$ws.Range("E5").Value = 14.69890766435346
$ws.Range("G5").Value = 27.0131053990796
$ws.Range("H5").Value = 13.87576262031505
$ws.Range("K5").Value = 7.956771180255407
$ws.Range("O5").Value = 20.88377568512181

$ws.Range("B6").Value = 8.604798660941466
$ws.Range("C6").Value = 6.137112197867718
$ws.Range("D6").Value = 5.748877880388744
$ws.Range("E6").Value = 14.65713949843556
$ws.Range("G6").Value = 27.01497828112703
$ws.Range("H6").Value = 13.87811924312698
$ws.Range("K6").Value = 7.941267392227671
$ws.Range("O6").Value = 20.88726957906349

$ws.Range("B7").Value = 8.714417156320014
$ws.Range("C7").Value = 6.183831821983091
$ws.Range("D7").Value = 5.783129351534159
$ws.Range("E7").Value = 14.94499645996254
$ws.Range("G7").Value = 27.00293433068208
$ws.Range("H7").Value = 13.8619936239859
$ws.Range("K7").Value = 8.048095247610384
$ws.Range("O7").Value = 20.86349753190052

$ws.Range("B8").Value = 9.184808708220974
$ws.Range("C8").Value = 6.384253278870827
$ws.Range("D8").Value = 5.935659802405695
$ws.Range("E8").Value = 16.1622384628353
$ws.Range("G8").Value = 26.97353346118482
$ws.Range("H8").Value = 13.79658077429199
$ws.Range("K8").Value = 8.49950344135466
$ws.Range("O8").Value = 20.77049839942238

$ws.Range("B9").Value = 10.05072072456783
$ws.Range("C9").Value = 6.754327402936585
$ws.Range("D9").Value = 6.236721745447356
$ws.Range("E9").Value = 18.42313885602556
$ws.Range("G9").Value = 27.00191904782558
$ws.Range("H9").Value = 13.68902670921545
$ws.Range("K9").Value = 9.308296222581507
$ws.Range("O9").Value = 20.63126418429509

$ws.Range("B10").Value = 10.64611861107303
$ws.Range("C10").Value = 7.010389500585874
$ws.Range("D10").Value = 6.45649531895684
$ws.Range("E10").Value = 20.05936710415525
$ws.Range("G10").Value = 27.07563116908796
$ws.Range("H10").Value = 13.62270071860918
$ws.Range("K10").Value = 9.852736994862139
$ws.Range("O10").Value = 20.55558022621615

$ws.Range("B11").Value = 10.90705355415064
$ws.Range("C11").Value = 7.123100786117787
$ws.Range("D11").Value = 6.555650169563235
$ws.Range("E11").Value = 20.76158916221022
$ws.Range("G11").Value = 27.12067051463908
$ws.Range("H11").Value = 13.59529701781178
$ws.Range("K11").Value = 10.08909105787929
$ws.Range("O11").Value = 20.52699042379795

$ws.Range("B12").Value = 11.004363622593
$ws.Range("C12").Value = 7.165213814046073
$ws.Range("D12").Value = 6.593039079519052
$ws.Range("E12").Value = 21.0214836220469
$ws.Range("G12").Value = 27.13937816185733
$ws.Range("H12").Value = 13.58531915505154
$ws.Range("K12").Value = 10.176930578445
$ws.Range("O12").Value = 20.51700857077698

$ws.Range("B13").Value = 10.9834739643061
$ws.Range("C13").Value = 7.156169681129229
$ws.Range("D13").Value = 6.584994399540394
$ws.Range("E13").Value = 20.96577799700194
$ws.Range("G13").Value = 27.13527572315926
$ws.Range("H13").Value = 13.58745028659876
$ws.Range("K13").Value = 10.15808721831715
$ws.Range("O13").Value = 20.51912070924969

$ws.Range("B14").Value = 10.91508977152823
$ws.Range("C14").Value = 7.126576979217365
$ws.Range("D14").Value = 6.558729587188957
$ws.Range("E14").Value = 20.78309133304261
$ws.Range("G14").Value = 27.1221765157485
$ws.Range("H14").Value = 13.59446812291079
$ws.Range("K14").Value = 10.09635117330652
$ws.Range("O14").Value = 20.52615225949946

$ws.Range("B15").Value = 10.87300505430957
$ws.Range("C15").Value = 7.108375860506203
$ws.Range("D15").Value = 6.542619769705979
$ws.Range("E15").Value = 20.67040723816322
$ws.Range("G15").Value = 27.11436791704779
$ws.Range("H15").Value = 13.5988187907575
$ws.Range("K15").Value = 10.05831858090199
$ws.Range("O15").Value = 20.53056940168279

$ws.Range("B16").Value = 10.62885977213842
$ws.Range("C16").Value = 7.002945380717081
$ws.Range("D16").Value = 6.449995251397047
$ws.Range("E16").Value = 20.01263078392628
$ws.Range("G16").Value = 27.07291915027705
$ws.Range("H16").Value = 13.62454742708324
$ws.Range("K16").Value = 9.837059844283988
$ws.Range("O16").Value = 20.55756655504701

$ws.Range("B17").Value = 10.47648870778203
$ws.Range("C17").Value = 6.937281848895182
$ws.Range("D17").Value = 6.392933020017984
$ws.Range("E17").Value = 19.5983453721694
$ws.Range("G17").Value = 27.05043811256483
$ws.Range("H17").Value = 13.64104103340427
$ws.Range("K17").Value = 9.698399602808754
$ws.Range("O17").Value = 20.57562751355466

$ws.Range("B18").Value = 10.38791944803208
$ws.Range("C18").Value = 6.899160215197685
$ws.Range("D18").Value = 6.36003733152292
$ws.Range("E18").Value = 19.35609414795806
$ws.Range("G18").Value = 27.03859107423427
$ws.Range("H18").Value = 13.6507882135337
$ws.Range("K18").Value = 9.617582629885856
$ws.Range("O18").Value = 20.58656504333556

$ws.Range("B19").Value = 10.35777423191837
$ws.Range("C19").Value = 6.886192959079067
$ws.Range("D19").Value = 6.34888781538533
$ws.Range("E19").Value = 19.27339000577048
$ws.Range("G19").Value = 27.03476598497091
$ws.Range("H19").Value = 13.65413314337649
$ws.Range("K19").Value = 9.590037954013935
$ws.Range("O19").Value = 20.59036250197343

$ws.Range("B20").Value = 10.49280563657025
$ws.Range("C20").Value = 6.944308645397341
$ws.Range("D20").Value = 6.399015457223022
$ws.Range("E20").Value = 19.64285690601157
$ws.Range("G20").Value = 27.05271913226624
$ws.Range("H20").Value = 13.63925829164356
$ws.Range("K20").Value = 9.713270542693056
$ws.Range("O20").Value = 20.57364800277685

$ws.Range("B21").Value = 10.93521712850574
$ws.Range("C21").Value = 7.135284690496449
$ws.Range("D21").Value = 6.566448824660472
$ws.Range("E21").Value = 20.83691400210287
$ws.Range("G21").Value = 27.12597926406239
$ws.Range("H21").Value = 13.5923959668791
$ws.Range("K21").Value = 10.11452992493823
$ws.Range("O21").Value = 20.52406396678117

$ws.Range("B22").Value = 11.21558601356134
$ws.Range("C22").Value = 7.256775781376279
$ws.Range("D22").Value = 6.674932888101377
$ws.Range("E22").Value = 21.58221489536689
$ws.Range("G22").Value = 27.18348620662922
$ws.Range("H22").Value = 13.5640968636457
$ws.Range("K22").Value = 10.36707240268514
$ws.Range("O22").Value = 20.49658225610581

$ws.Range("B23").Value = 11.06677257985109
$ws.Range("C23").Value = 7.192245668045673
$ws.Range("D23").Value = 6.617132081203756
$ws.Range("E23").Value = 21.18763315012625
$ws.Range("G23").Value = 27.1519143608328
$ws.Range("H23").Value = 13.57898719533599
$ws.Range("K23").Value = 10.23318393373265
$ws.Range("O23").Value = 20.51079767368918

$ws.Range("B24").Value = 10.4854317658321
$ws.Range("C24").Value = 6.941132983710466
$ws.Range("D24").Value = 6.396265864120683
$ws.Range("E24").Value = 19.62274593123593
$ws.Range("G24").Value = 27.05168452630263
$ws.Range("H24").Value = 13.64006344444479
$ws.Range("K24").Value = 9.706550812321904
$ws.Range("O24").Value = 20.57454121406955

$ws.Range("B25").Value = 9.823226913256015
$ws.Range("C25").Value = 6.656869107142467
$ws.Range("D25").Value = 6.155342833381821
$ws.Range("E25").Value = 17.78302421084269
$ws.Range("G25").Value = 26.98497612115555
$ws.Range("H25").Value = 13.71589822598354
$ws.Range("K25").Value = 9.098038967930682
$ws.Range("O25").Value = 20.66427926304305
